$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove data rows 3-8, keeping only the single remaining data row (row 2)
$ws.Range("A3:A8").EntireRow.Delete()

# Remove column E (the duplicate "rejection-f" / old E column), shifting D left is not needed
# since the new layout only spans A:D
$ws.Range("E1").EntireColumn.Delete()

# Header row: column C becomes "prediction", column D becomes "rejection-f"
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"

# Remaining data row now describes RUG520.fasta with updated values
$ws.Range("A2").Value = "RUG520.fasta"
$ws.Range("B2").Value = 653556.7941053929
$ws.Range("C2").Value = "o__Chitinivibrionales"
$ws.Range("D2").Value = "o__Chitinivibrionales"
